# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.690.84"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "2.768.44"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.00%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "578.64"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "161.19"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  -0.12%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.604"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "

$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("E10").Value = "  +4.85%  "

$ws.Range("E11").Value = "  +3.16%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.389"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "3.254.01"
$ws.Range("E13").Value = "  +0.48%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.43"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").Value = "64.213.35"
$ws.Range("E15").Value = "  +0.65%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0000153"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "2.770.14"
$ws.Range("E17").Value = "  +0.62%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "12.22"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.87"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "359.71"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.71"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("E22").Value = "  -0.07%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.532"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -6.60%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "65.31"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("E25").Value = "  -0.82%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.66"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").Value = "0.0₃0932"
$ws.Range("E28").Value = "  -0.77%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.40"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.90%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.39"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +9.94%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.12%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "168.09"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.03"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.53"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +3.53%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "20.24"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.85"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("E38").Value = "  -0.82%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "354.21"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +6.38%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.44"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.55%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "4.21"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.51%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "39.13"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.11%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "22.67"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "21.66"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("E45").Value = "  -0.62%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "137.17"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.633"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("E48").Value = "  -1.92%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.102"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").Value = "2.155.61"
$ws.Range("E50").Value = "  +1.33%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "

